# "Correcao dos 22 e 23" - rebuild the traceability matrix (Matriz de Rastreabilidade)
# with the corrected/extended list of system features (rows 2-26) and make room
# for a few extra blank rows (27-33) at the end, matching the author's re-edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-26: "#", "Caracteristicas", "Requisitos" --------------------------
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 'Sistema Web '
$ws.Range("C2").Value = 'SSS15,SSS13 ,SSS08'

$ws.Range("A3").Value = 14
$ws.Range("B3").Value = 'Gestão de perfil do usuário '
$ws.Range("C3").Value = 'SSS15'

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 'Registro de venda '
$ws.Range("C4").Value = 'SSS03,SSS04,SSS05,SSS12 '

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 'Realização  de orçamento '
$ws.Range("C5").Value = 'SSS17'

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 'Recebimento de  mercadorias '
$ws.Range("C6").Value = 'SSS01,SSS02'

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 'Contagem do estoque  '
$ws.Range("C7").Value = 'SSS02'

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 'Vendas de produtos '
$ws.Range("C8").Value = 'SSS18'

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = 'Trocas produtos '
$ws.Range("C9").Value = 'SSS10'

$ws.Range("A10").Value = 25
$ws.Range("B10").Value = 'Devolução de mercadoria '
$ws.Range("C10").Value = 'SSS19'

$ws.Range("A11").Value = 20
$ws.Range("B11").Value = 'Catalogo de produtos '
$ws.Range("C11").Value = 'SSS21'

$ws.Range("A12").Value = 18
$ws.Range("B12").Value = 'Avisos via telefone '
$ws.Range("C12").Value = 'SSS22'

$ws.Range("A13").Value = 19
$ws.Range("B13").Value = 'Divulgação da marca '
$ws.Range("C13").Value = 'SSS21'

$ws.Range("A14").Value = 22
$ws.Range("B14").Value = 'Divulgação de produtos '
$ws.Range("C14").Value = 'SSS21'

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = 'Notificação de venda  '
$ws.Range("C15").Value = 'SSS23'

$ws.Range("A16").Value = 16
$ws.Range("B16").Value = 'Aprovação de orçamento do material '
$ws.Range("C16").Value = 'SSS24'

$ws.Range("A17").Value = 8
$ws.Range("B17").Value = 'Consulta de preço '
$ws.Range("C17").Value = 'SSS25'

$ws.Range("A18").Value = 24
$ws.Range("B18").Value = 'Alteração de prazo  '
$ws.Range("C18").Value = 'SSS13,SSS17'

$ws.Range("A19").Value = 6
$ws.Range("B19").Value = 'Verificar material do produto '
$ws.Range("C19").Value = 'SSS24'

$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 'Relatório de venda  '
$ws.Range("C20").Value = 'SSS11 '

$ws.Range("A21").Value = 9
$ws.Range("B21").Value = 'Inventário de produtos'
$ws.Range("C21").Value = 'SSS20'

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 'Alteração de preço '
$ws.Range("C22").Value = 'SSS11 '

$ws.Range("A23").Value = 23
$ws.Range("B23").Value = 'Promoção de vendas '
$ws.Range("C23").Value = 'SSS21'

$ws.Range("A24").Value = 13
$ws.Range("B24").Value = 'Avisos via email '
$ws.Range("C24").Value = 'SSS07'

$ws.Range("A25").Value = 15
$ws.Range("B25").Value = 'Avisos via whatsapp '
$ws.Range("C25").Value = 'SSS07'

$ws.Range("A26").Value = 17
$ws.Range("B26").Value = 'Avisos via facebook '
$ws.Range("C26").Value = 'SSS07'

# New rows 9-26 need column A centered (same look as the rest of column A)
$ws.Range("A9:A26").VerticalAlignment = -4108
$ws.Range("A9:A26").HorizontalAlignment = -4108

# --- Rows 27-33: trailing blank rows (column A only), same centered style ----
for ($r = 27; $r -le 33; $r++) {
    $ws.Range("A$r").VerticalAlignment = -4108
    $ws.Range("A$r").HorizontalAlignment = -4108
}

# Grow the "Tabela1" table/autofilter to cover the new data range
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C26"))

# Leave the selection where the author left it
$ws.Range("C19").Select() | Out-Null
